# Add two new columns (I: "I0", J: "IF") to the single worksheet.
# Header row (row 1) gets the same style as the existing header cells (H1 etc.),
# and data rows 2-38 get the new numeric values, mirroring the existing
# unstyled data cells in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# I0 / IF values per row (2..38), taken straight from the new data.
$values = @{
    2  = @(4, 5)
    3  = @(10, 10)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(5, 9)
    7  = @(1, 4)
    8  = @(1, 7)
    9  = @(1, 7)
    10 = @(1, 7)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 5)
    15 = @(1, 5)
    16 = @(1, 6)
    17 = @(1, 5)
    18 = @(1, 7)
    19 = @(1, 6)
    20 = @(1, 7)
    21 = @(1, 6)
    22 = @(1, 5)
    23 = @(1, 6)
    24 = @(1, 5)
    25 = @(1, 7)
    26 = @(1, 6)
    27 = @(1, 7)
    28 = @(1, 6)
    29 = @(1, 7)
    30 = @(1, 7)
    31 = @(1, 5)
    32 = @(1, 8)
    33 = @(1, 6)
    34 = @(1, 5)
    35 = @(1, 5)
    36 = @(1, 4)
    37 = @(1, 4)
    38 = @(1, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
